$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep text formatting so values are not auto-converted to numbers
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "B8", "C8", "D8", "E8", "B9", "C9", "D9", "E9", "B10", "C10", "D10", "E10", "B11", "C11", "D11", "E11", "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "E49", "E50", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "319.40"
$ws.Range("E2").Value = "4.55%"
$ws.Range("D3").Value = "35.99"
$ws.Range("E3").Value = "-0.22%"
$ws.Range("D4").Value = "5.120"
$ws.Range("E4").Value = "0.78%"
$ws.Range("D5").Value = "0.08204"
$ws.Range("E5").Value = "4.44%"
$ws.Range("D6").Value = "2.145"
$ws.Range("E6").Value = "-1.31%"
$ws.Range("D7").Value = "8.038"
$ws.Range("E7").Value = "1.48%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.139"
$ws.Range("E8").Value = "1.04%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9256"
$ws.Range("E9").Value = "0.36%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1006"
$ws.Range("E10").Value = "4.11%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1900"
$ws.Range("E11").Value = "1.80%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09195"
$ws.Range("E12").Value = "5.56%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03611"
$ws.Range("E13").Value = "3.54%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09922"
$ws.Range("E14").Value = "0.01%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001439"
$ws.Range("E15").Value = "0.53%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005728"
$ws.Range("E16").Value = "1.32%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").Value = "0.12%"
$ws.Range("D18").Value = "2.798"
$ws.Range("E18").Value = "16.75%"
$ws.Range("D19").Value = "0.3373"
$ws.Range("E19").Value = "-1.50%"
$ws.Range("E20").Value = "2.35%"
$ws.Range("D21").Value = "5.086"
$ws.Range("E21").Value = "4.72%"
$ws.Range("E22").Value = "-0.51%"
$ws.Range("D23").Value = "0.04604"
$ws.Range("E23").Value = "1.24%"
$ws.Range("D24").Value = "0.001244"
$ws.Range("E24").Value = "1.08%"
$ws.Range("E25").Value = "-6.95%"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "-7.11%"
$ws.Range("E27").Value = "-5.23%"
$ws.Range("D39").Value = "0.02008"
$ws.Range("E39").Value = "9.65%"
$ws.Range("D40").Value = "0.04987"
$ws.Range("E40").Value = "4.37%"
$ws.Range("D41").Value = "0.007793"
$ws.Range("E41").Value = "1.26%"
$ws.Range("D42").Value = "0.1399"
$ws.Range("E42").Value = "0.09%"
$ws.Range("D43").Value = "0.007831"
$ws.Range("E43").Value = "1.33%"
$ws.Range("D44").Value = "0.002131"
$ws.Range("E44").Value = "-4.75%"
$ws.Range("D45").Value = "0.01203"
$ws.Range("E45").Value = "7.22%"
$ws.Range("D46").Value = "0.00006444"
$ws.Range("E46").Value = "1.14%"
$ws.Range("E47").Value = "0.04%"
$ws.Range("E48").Value = "17.90%"
$ws.Range("E49").Value = "-4.96%"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").Value = "0.04%"
